# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off again (new xliff files generated) and is no longer in sync,
# while a.md stays untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69721c65272199cff0f741c297ff9e49804c8bc8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9839aa20c1325119189fc11266d61d3298b02a95/e2e/b.md."

# ---- Overview sheet: row 3 is b.md ----
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-01 08:47:52"

# ---- zh-cn sheet: row 3 is b.md ----
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-01 08:47:47"
$zhcn.Range("P3").Value = $errorDetail

# ---- de-de sheet: row 3 is b.md ----
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-01 08:47:52"
$dede.Range("P3").Value = $errorDetail

# The new "Error Detail" text is much longer than the prior placeholder,
# so the column is widened to fit the report (ColumnWidth 39.17 stores as
# the OOXML column width of 40 characters).
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
